$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4149273333333334
$ws.Range("H2").Value = 1.244782
$ws.Range("I2").Value = 0.1353844755004719
$ws.Range("J2").Value = 0.1353844755004719
$ws.Range("M2").Value = 19.59940166666667
$ws.Range("N2").Value = 58.798205
$ws.Range("O2").Value = 0.1807871245579405
$ws.Range("P2").Value = 0.1807871245579405
$ws.Range("Q2").Value = 8.13232746847889
$ws.Range("R2").Value = 73.19094721630999
$ws.Range("S2").Value = 0.02447577003551525
$ws.Range("T2").Value = 0.02447577003551525
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4149273333333334
$ws.Range("H3").Value = 1.244782
$ws.Range("I3").Value = 0.1353844755004719
$ws.Range("J3").Value = 0.1353844755004719
$ws.Range("O3").Value = 0.1000607063571047
$ws.Range("P3").Value = 0.1000607063571047
$ws.Range("Q3").Value = 4.501019819929112
$ws.Range("R3").Value = 40.50917837936201
$ws.Range("S3").Value = 0.01354666624836336
$ws.Range("T3").Value = 0.01354666624836336
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4149273333333334
$ws.Range("H4").Value = 1.244782
$ws.Range("I4").Value = 0.1353844755004719
$ws.Range("J4").Value = 0.1353844755004719
$ws.Range("M4").Value = 4.842319
$ws.Range("N4").Value = 14.526957
$ws.Range("O4").Value = 0.04466610476640988
$ws.Range("P4").Value = 0.04466610476640988
$ws.Range("Q4").Value = 2.009210509819333
$ws.Range("R4").Value = 18.082894588374
$ws.Range("S4").Value = 0.00604709716644953
$ws.Range("T4").Value = 0.00604709716644953
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4149273333333334
$ws.Range("H5").Value = 1.244782
$ws.Range("I5").Value = 0.1353844755004719
$ws.Range("J5").Value = 0.1353844755004719
$ws.Range("M5").Value = 73.12203966666667
$ws.Range("N5").Value = 219.366119
$ws.Range("O5").Value = 0.674486064318545
$ws.Range("P5").Value = 0.674486064318545
$ws.Range("Q5").Value = 30.34033292678422
$ws.Range("R5").Value = 273.062996341058
$ws.Range("S5").Value = 0.09131494205014377
$ws.Range("T5").Value = 0.09131494205014377
$ws.Range("I6").Value = 0.389303862711544
$ws.Range("J6").Value = 0.389303862711544
$ws.Range("M6").Value = 19.59940166666667
$ws.Range("N6").Value = 58.798205
$ws.Range("O6").Value = 0.1807871245579405
$ws.Range("P6").Value = 0.1807871245579405
$ws.Range("Q6").Value = 23.38485623710222
$ws.Range("R6").Value = 210.46370613392
$ws.Range("S6").Value = 0.07038112591891926
$ws.Range("T6").Value = 0.07038112591891926
$ws.Range("I7").Value = 0.389303862711544
$ws.Range("J7").Value = 0.389303862711544
$ws.Range("O7").Value = 0.1000607063571047
$ws.Range("P7").Value = 0.1000607063571047
$ws.Range("S7").Value = 0.03895401949046642
$ws.Range("T7").Value = 0.03895401949046642
$ws.Range("I8").Value = 0.389303862711544
$ws.Range("J8").Value = 0.389303862711544
$ws.Range("M8").Value = 4.842319
$ws.Range("N8").Value = 14.526957
$ws.Range("O8").Value = 0.04466610476640988
$ws.Range("P8").Value = 0.04466610476640988
$ws.Range("Q8").Value = 5.777570948085334
$ws.Range("R8").Value = 51.998138532768
$ws.Range("S8").Value = 0.01738868711784188
$ws.Range("T8").Value = 0.01738868711784188
$ws.Range("I9").Value = 0.389303862711544
$ws.Range("J9").Value = 0.389303862711544
$ws.Range("M9").Value = 73.12203966666667
$ws.Range("N9").Value = 219.366119
$ws.Range("O9").Value = 0.674486064318545
$ws.Range("P9").Value = 0.674486064318545
$ws.Range("Q9").Value = 87.24492790393955
$ws.Range("R9").Value = 785.204351135456
$ws.Range("S9").Value = 0.2625800301843165
$ws.Range("T9").Value = 0.2625800301843165
$ws.Range("G10").Value = 1.270157666666667
$ws.Range("H10").Value = 3.810473
$ws.Range("I10").Value = 0.4144331204288861
$ws.Range("J10").Value = 0.4144331204288861
$ws.Range("M10").Value = 19.59940166666667
$ws.Range("N10").Value = 58.798205
$ws.Range("O10").Value = 0.1807871245579405
$ws.Range("P10").Value = 0.1807871245579405
$ws.Range("Q10").Value = 24.89433028899611
$ws.Range("R10").Value = 224.048972600965
$ws.Range("S10").Value = 0.07492417216391298
$ws.Range("T10").Value = 0.07492417216391298
$ws.Range("G11").Value = 1.270157666666667
$ws.Range("H11").Value = 3.810473
$ws.Range("I11").Value = 0.4144331204288861
$ws.Range("J11").Value = 0.4144331204288861
$ws.Range("O11").Value = 0.1000607063571047
$ws.Range("P11").Value = 0.1000607063571047
$ws.Range("Q11").Value = 13.77832784881589
$ws.Range("R11").Value = 124.004950639343
$ws.Range("S11").Value = 0.0414684707678934
$ws.Range("T11").Value = 0.0414684707678934
$ws.Range("G12").Value = 1.270157666666667
$ws.Range("H12").Value = 3.810473
$ws.Range("I12").Value = 0.4144331204288861
$ws.Range("J12").Value = 0.4144331204288861
$ws.Range("M12").Value = 4.842319
$ws.Range("N12").Value = 14.526957
$ws.Range("O12").Value = 0.04466610476640988
$ws.Range("P12").Value = 0.04466610476640988
$ws.Range("Q12").Value = 6.150508602295667
$ws.Range("R12").Value = 55.35457742066099
$ws.Range("S12").Value = 0.01851111317574679
$ws.Range("T12").Value = 0.01851111317574679
$ws.Range("G13").Value = 1.270157666666667
$ws.Range("H13").Value = 3.810473
$ws.Range("I13").Value = 0.4144331204288861
$ws.Range("J13").Value = 0.4144331204288861
$ws.Range("M13").Value = 73.12203966666667
$ws.Range("N13").Value = 219.366119
$ws.Range("O13").Value = 0.674486064318545
$ws.Range("P13").Value = 0.674486064318545
$ws.Range("Q13").Value = 92.87651928492079
$ws.Range("R13").Value = 835.888673564287
$ws.Range("S13").Value = 0.279529364321333
$ws.Range("T13").Value = 0.279529364321333
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.186581
$ws.Range("H14").Value = 0.559743
$ws.Range("I14").Value = 0.06087854135909794
$ws.Range("J14").Value = 0.06087854135909794
$ws.Range("M14").Value = 19.59940166666667
$ws.Range("N14").Value = 58.798205
$ws.Range("O14").Value = 0.1807871245579405
$ws.Range("P14").Value = 0.1807871245579405
$ws.Range("Q14").Value = 3.656875962368333
$ws.Range("R14").Value = 32.911883661315
$ws.Range("S14").Value = 0.01100605643959297
$ws.Range("T14").Value = 0.01100605643959297
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.186581
$ws.Range("H15").Value = 0.559743
$ws.Range("I15").Value = 0.06087854135909794
$ws.Range("J15").Value = 0.06087854135909794
$ws.Range("O15").Value = 0.1000607063571047
$ws.Range("P15").Value = 0.1000607063571047
$ws.Range("Q15").Value = 2.023980373323667
$ws.Range("R15").Value = 18.215823359913
$ws.Range("S15").Value = 0.006091549850381554
$ws.Range("T15").Value = 0.006091549850381555
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.186581
$ws.Range("H16").Value = 0.559743
$ws.Range("I16").Value = 0.06087854135909794
$ws.Range("J16").Value = 0.06087854135909794
$ws.Range("M16").Value = 4.842319
$ws.Range("N16").Value = 14.526957
$ws.Range("O16").Value = 0.04466610476640988
$ws.Range("P16").Value = 0.04466610476640988
$ws.Range("Q16").Value = 0.9034847213389999
$ws.Range("R16").Value = 8.131362492051
$ws.Range("S16").Value = 0.002719207306371685
$ws.Range("T16").Value = 0.002719207306371685
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.186581
$ws.Range("H17").Value = 0.559743
$ws.Range("I17").Value = 0.06087854135909794
$ws.Range("J17").Value = 0.06087854135909794
$ws.Range("M17").Value = 73.12203966666667
$ws.Range("N17").Value = 219.366119
$ws.Range("O17").Value = 0.674486064318545
$ws.Range("P17").Value = 0.674486064318545
$ws.Range("Q17").Value = 13.64318328304633
$ws.Range("R17").Value = 122.788649547417
$ws.Range("S17").Value = 0.04106172776275174
$ws.Range("T17").Value = 0.04106172776275174
